# Fix typo noted while teaching the class: color the `http` inside
# `'http'` red on the "Creating an HTTPS Client and Server" slide
# (Slide 16, Content Placeholder 2, first paragraph).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange

# Locate "'http'" within the paragraph text so we don't depend on a
# hard-coded character offset.
$fullText = $tr.Text
$needle = "'http'"
$idx = $fullText.IndexOf($needle)

if ($idx -ge 0) {
    # PowerPoint TextRange.Characters is 1-indexed.
    $start = $idx + 1
    $httpOnly = $tr.Characters($start + 1, 4)
    $httpOnly.Font.Color.RGB = 255
}
